$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("checkBoxClickData")
$ws.Range("A3").Value = "Failed for Raport Purpoose"
$ws.Range("A3").Select()
